# Remove the "productgroups" row from the Table2 table on the "Tables" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")

# Delete the entire worksheet row 12 ("productgroups"), which is part of
# Table2 (A1:E27). This shifts subsequent rows up and shrinks the table
# range to A1:E26.
$ws.Rows("12:12").Delete()

# Restore the active selection to match the post-edit workbook state.
$ws.Range("B14").Select()
